# Apply the register-list template changes to the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Remove the unused Sheet2 / Sheet3 worksheets, keeping only Sheet1.
# ---------------------------------------------------------------------
foreach ($name in @("Sheet2", "Sheet3")) {
    $sheet = $wb.Worksheets.Item($name)
    $sheet.Delete()
}

$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 2. Move the RW/R access-type value from the first field row of each
#    5-row register block up onto the block's header row, leaving the
#    field rows blank in column E.
# ---------------------------------------------------------------------
$headerRows = 7..52 | Where-Object { (($_ - 7) % 5) -eq 0 }

foreach ($h in $headerRows) {
    $fieldRow = $h + 1
    $val = $ws.Cells.Item($fieldRow, 5).Value2
    $ws.Cells.Item($h, 5).Value = $val
    $ws.Cells.Item($fieldRow, 5).Value = $null
    $ws.Cells.Item($fieldRow + 1, 5).Value = $null
    $ws.Cells.Item($fieldRow + 2, 5).Value = $null
    $ws.Cells.Item($fieldRow + 3, 5).Value = $null
}

# ---------------------------------------------------------------------
# 3. Update the active view (scroll position + selection).
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("I43").Select()
